# Append 12 new daily COVID-19 data rows (2020-05-25 .. 2020-06-05) to the
# "Covid-19 podatki" sheet, growing the Tabela1 table from A1:J75 to A1:J87.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, Date(serial), Tested(all), Tested(daily), Positive(all),
#             Positive(daily), Hospitalized, ICU, Discharged, Deaths(all), Deaths(daily)
$rows = @(
    @(76, 43976, 75770, 754, 1469, 0, 9, 2, 6, 108, 1),
    @(77, 43977, 76579, 809, 1471, 2, 8, 2, 2, 108, 0),
    @(78, 43978, 77210, 631, 1473, 2, 7, 2, 1, 108, 0),
    @(79, 43979, 77916, 706, 1473, 0, 7, 2, 0, 108, 0),
    @(80, 43980, 78529, 613, 1473, 0, 7, 2, 0, 108, 0),
    @(81, 43981, 78793, 264, 1473, 0, 6, 2, 1, 108, 0),
    @(82, 43982, 79039, 246, 1473, 0, 5, 1, 0, 109, 1),
    @(83, 43983, 79698, 659, 1475, 2, 5, 1, 0, 109, 0),
    @(84, 43984, 80505, 807, 1477, 2, 5, 0, 0, 109, 0),
    @(85, 43985, 81333, 828, 1477, 0, 5, 0, 0, 109, 0),
    @(86, 43986, 82161, 828, 1479, 2, 6, 0, 0, 109, 0),
    @(87, 43987, 82876, 715, 1484, 5, 6, 0, 0, 109, 0)
)

$cols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")

foreach ($entry in $rows) {
    $r = $entry[0]

    # Push a fresh row into place, inheriting the cell formatting of the row
    # directly above it (matches how Excel extends a banded table when a new
    # row is typed in underneath it).
    $ws.Rows([string]$r).Insert(-4121, 0)

    for ($i = 1; $i -le 10; $i++) {
        $addr = $cols[$i - 1] + [string]$r
        $ws.Range($addr).Value = [double]$entry[$i]
    }
}

# Grow the table / autofilter to cover the newly added rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:J87"))

# Match the final view/selection state saved in the workbook.
$win = $excel.ActiveWindow
$win.ScrollRow = 65
$win.ScrollColumn = 1
$ws.Range("A87:J87").Select()
